# Apply the crypto price/volume refresh described in the commit diff.
# Column order per row: B=Coin, C=Link, D=Price, E=Volume(1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.263.98'
$ws.Range('E2').Value = '  -0.96%  '
$ws.Range('D3').Value = '1.783.53'
$ws.Range('E3').Value = '  -2.18%  '
$ws.Range('D4').Value = '''1.005'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '''340.27'
$ws.Range('E5').Value = '  -0.81%  '
$ws.Range('D6').Value = '''1.001'
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').Value = '''0.3990'
$ws.Range('E7').Value = '  +4.50%  '
$ws.Range('D8').Value = '''0.3445'
$ws.Range('E8').Value = '  -2.53%  '
$ws.Range('D9').Value = '''47.89'
$ws.Range('E9').Value = '  -4.01%  '
$ws.Range('D10').Value = '''1.192'
$ws.Range('E10').Value = '  -3.67%  '
$ws.Range('D11').Value = '''0.07446'
$ws.Range('E11').Value = '  -3.72%  '
$ws.Range('D12').Value = '''1.002'
$ws.Range('E12').Value = '  +0.03%  '
$ws.Range('D13').Value = '''21.66'
$ws.Range('E13').Value = '  -2.35%  '
$ws.Range('D14').Value = '''6.458'
$ws.Range('E14').Value = '  -2.31%  '
$ws.Range('D15').Value = '1.785.41'
$ws.Range('E15').Value = '  -2.36%  '
$ws.Range('D16').Value = '''7.101'
$ws.Range('E16').Value = '  -1.59%  '
$ws.Range('D17').Value = '''0.00001089'
$ws.Range('E17').Value = '  -3.22%  '
$ws.Range('D18').Value = '''0.06688'
$ws.Range('E18').Value = '  -0.77%  '
$ws.Range('D19').Value = '''84.07'
$ws.Range('E19').Value = '  -3.42%  '
$ws.Range('E20').Value = '  -0.19%  '
$ws.Range('D21').Value = '''17.62'
$ws.Range('E21').Value = '  +0.21%  '
$ws.Range('D22').Value = '''6.491'
$ws.Range('E22').Value = '  -0.62%  '
$ws.Range('D23').Value = '27.279.11'
$ws.Range('E23').Value = '  -0.88%  '
$ws.Range('D24').Value = '''12.37'
$ws.Range('E24').Value = '  -6.04%  '
$ws.Range('E25').Value = '  -3.85%  '
$ws.Range('D26').Value = '''1.468'
$ws.Range('E26').Value = '  -0.92%  '
$ws.Range('D27').Value = '''21.14'
$ws.Range('E27').Value = '  -3.94%  '
$ws.Range('D28').Value = '''2.486'
$ws.Range('E28').Value = '  -7.47%  '
$ws.Range('D29').Value = '''157.47'
$ws.Range('E29').Value = '  +2.94%  '
$ws.Range('D30').Value = '1.987.20'
$ws.Range('E30').Value = '  -2.35%  '
$ws.Range('D31').Value = '''135.70'
$ws.Range('E31').Value = '  +0.20%  '
$ws.Range('D32').Value = '''4.024'
$ws.Range('E32').Value = '  -1.61%  '
$ws.Range('D33').Value = '''5.942'
$ws.Range('E33').Value = '  -6.43%  '
$ws.Range('D34').Value = '''0.08853'
$ws.Range('E34').Value = '  +0.70%  '
$ws.Range('D35').Value = '''12.94'
$ws.Range('E35').Value = '  -6.89%  '
# Row 36: coin replaced/reordered
$ws.Range('B36').Value = 'WEMIXTOKEN'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').Value = '''1.621'
$ws.Range('E36').Value = '  -4.46%  '
# Row 37: coin replaced/reordered
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '''0.02424'
$ws.Range('E37').Value = '  +0.70%  '
$ws.Range('D38').Value = '''5.392'
$ws.Range('E38').Value = '  -4.07%  '
# Row 39: coin replaced/reordered
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '''0.06446'
$ws.Range('E39').Value = '  -1.14%  '
# Row 40: coin replaced/reordered
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '''0.6796'
$ws.Range('E40').Value = '  -2.95%  '
$ws.Range('D41').Value = '''0.2196'
$ws.Range('E41').Value = '  -2.84%  '
$ws.Range('D42').Value = '''1.254'
$ws.Range('E42').Value = '  -3.84%  '
$ws.Range('D43').Value = '''8.366'
$ws.Range('E43').Value = '  -8.21%  '
$ws.Range('D44').Value = '''14.43'
$ws.Range('E44').Value = '  -1.72%  '
$ws.Range('D45').Value = '''1.000'
$ws.Range('E45').Value = '  -0.16%  '
$ws.Range('D46').Value = '''0.6377'
$ws.Range('E46').Value = '  -3.63%  '
$ws.Range('D47').Value = '''3.876'
$ws.Range('E47').Value = '  -1.71%  '
$ws.Range('D48').Value = '''132.40'
$ws.Range('E48').Value = '  -0.62%  '
$ws.Range('D49').Value = '''2.127'
$ws.Range('E49').Value = '  -2.85%  '
$ws.Range('E50').Value = '  -2.35%  '
# Row 51: coin replaced/reordered
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '''79.13'
$ws.Range('E51').Value = '  -2.42%  '
